$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Unprotect()
$ws.Range("A1").Value = "test"
$ws.Protect($null, $true, $true, $true, $true)
